$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New daily rows (date serial, nuovi pos., somma mobile 7gg., somma mobile 7gg. per 100mila abitanti)
$data = @(
    @(329, 44403, 0, 1, 9.930486593843098),
    @(330, 44404, 0, 1, 9.930486593843098),
    @(331, 44405, 0, 1, 9.930486593843098),
    @(332, 44406, 0, 1, 9.930486593843098),
    @(333, 44407, 1, 2, 19.8609731876862),
    @(334, 44408, 1, 3, 29.7914597815293),
    @(335, 44409, 1, 3, 29.7914597815293),
    @(336, 44410, 1, 4, 39.72194637537239),
    @(337, 44411, 0, 4, 39.72194637537239),
    @(338, 44412, 0, 4, 39.72194637537239),
    @(339, 44413, 3, 7, 69.51340615690168),
    @(340, 44414, 0, 6, 59.5829195630586),
    @(341, 44415, 0, 5, 49.65243296921549),
    @(342, 44416, 0, 4, 39.72194637537239),
    @(343, 44417, 4, 7, 69.51340615690168)
)

# Column A carries the date style (s="2") used throughout the column;
# copy that formatting (from the last existing data row) before writing
# the new serial-date values so the pasted rows match the existing ones.
$ws.Range("A328").Copy()

foreach ($row in $data) {
    $r = $row[0]
    $ws.Range("A$r").PasteSpecial(-4122)
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
}

$excel.CutCopyMode = 0
